$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 164 (shifts existing 164..172 down to 165..173,
# carrying their data/formatting with them automatically).
$ws.Rows.Item(164).Insert()

# Populate the freshly inserted row 164 with the new week's record.
$ws.Cells.Item(164, 1).Value = 11
$ws.Cells.Item(164, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(164, 3).Value = "Bíobío"
$ws.Cells.Item(164, 4).Value = 44568
$ws.Cells.Item(164, 5).Value = 8
$ws.Cells.Item(164, 6).Value = 100114013
$ws.Cells.Item(164, 7).Value = "Zanahoria"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 1200
$ws.Cells.Item(164, 11).Value = 7500
$ws.Cells.Item(164, 12).Value = 8500
$ws.Cells.Item(164, 13).Value = 8000
$ws.Cells.Item(164, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(164, 15).Value = "Región de Ñuble"
$ws.Cells.Item(164, 16).Value = 400
$ws.Cells.Item(164, 17).Value = 20
$ws.Cells.Item(164, 18).Value = "Hortaliza"
